$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 338.14285
$ws.Range("I2").Value = 327.83334
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 327.83334
$ws.Range("L2").Value = 400
$ws.Range("M2").Value = -214.83334
$ws.Range("N2").Value = -626

$ws.Range("H33").Value = 2274.8572
$ws.Range("I33").Value = 1889.4445
$ws.Range("J33").Value = 2968.6
$ws.Range("K33").Value = 1889.4445
$ws.Range("L33").Value = 2968.6
$ws.Range("M33").Value = -1660.4445
$ws.Range("N33").Value = -3426.6

$ws.Range("H53").Value = 616.61536
$ws.Range("I53").Value = 786.26666
$ws.Range("J53").Value = 385.27274
$ws.Range("K53").Value = 786.26666
$ws.Range("L53").Value = 385.27274
$ws.Range("M53").Value = -149.26666
$ws.Range("N53").Value = -1659.27274

$ws.Range("H74").Value = 6552.0527
$ws.Range("I74").Value = 6552.0527
$ws.Range("K74").Value = 6552.0527
$ws.Range("M74").Value = -5616.0527

$ws.Range("H77").Value = 6552.0527
$ws.Range("I77").Value = 6552.0527
$ws.Range("K77").Value = 32760.2635
$ws.Range("M77").Value = -28080.2635

$ws.Range("H80").Value = 1857.5
$ws.Range("I80").Value = 404.5
$ws.Range("J80").Value = 2148.1
$ws.Range("K80").Value = 1213.5
$ws.Range("L80").Value = 6444.299999999999
$ws.Range("M80").Value = -215.5
$ws.Range("N80").Value = -8440.299999999999

$ws.Range("H83").Value = 1857.5
$ws.Range("I83").Value = 404.5
$ws.Range("J83").Value = 2148.1
$ws.Range("K83").Value = 3640.5
$ws.Range("L83").Value = 19332.9
$ws.Range("M83").Value = 1351.5
$ws.Range("N83").Value = -29316.9

$ws.Range("H98").Value = 944.9
$ws.Range("I98").Value = 944.9
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 944.9
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("M98").Value = 553.1

$ws.Range("H111").Value = 3093.2942
$ws.Range("I111").Value = 3229.8462
$ws.Range("J111").Value = 2649.5
$ws.Range("K111").Value = 9689.5386
$ws.Range("L111").Value = 7948.5
$ws.Range("M111").Value = -6622.5386
$ws.Range("N111").Value = -14082.5

$ws.Range("H122").Value = 944.9
$ws.Range("I122").Value = 944.9
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2834.7
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("M122").Value = -384.6999999999998

$ws.Range("H135").Value = 1113.6666
$ws.Range("I135").Value = 1113.6666
$ws.Range("K135").Value = 10022.9994
$ws.Range("M135").Value = -7487.999400000001

$ws.Range("H141").Value = 3553.8667
$ws.Range("I141").Value = 3400.6365
$ws.Range("K141").Value = 10201.9095
$ws.Range("M141").Value = -5021.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2499
$ws.Range("I61").Value = 2499
$ws.Range("K61").Value = 2499
$ws.Range("M61").Value = -2287

$ws.Range("H74").Value = 1947.4445
$ws.Range("I74").Value = 1767.9412
$ws.Range("K74").Value = 1767.9412
$ws.Range("M74").Value = -893.9412

$ws.Range("H77").Value = 1947.4445
$ws.Range("I77").Value = 1767.9412
$ws.Range("K77").Value = 8839.706
$ws.Range("M77").Value = -4471.706

$ws.Range("H94").Value = 1000330
$ws.Range("J94").Value = 1000330
$ws.Range("L94").Value = 1000330
$ws.Range("N94").Value = -1002132

$ws.Range("H110").Value = 1366.3334
$ws.Range("I110").Value = 1366.3334
$ws.Range("K110").Value = 1366.3334
$ws.Range("M110").Value = 678.6666

$ws.Range("H136").Value = 2499
$ws.Range("I136").Value = 2499
$ws.Range("K136").Value = 7497
$ws.Range("M136").Value = -4947

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 70180
$ws.Range("J62").Value = 70180
$ws.Range("L62").Value = 70180
$ws.Range("N62").Value = -71552

$ws.Range("H65").Value = 70180
$ws.Range("J65").Value = 70180
$ws.Range("L65").Value = 210540
$ws.Range("N65").Value = -217404

$ws.Range("H86").Value = 12568.404
$ws.Range("I86").Value = 12579.56
$ws.Range("J86").Value = 12552
$ws.Range("K86").Value = 12579.56
$ws.Range("L86").Value = 12552
$ws.Range("M86").Value = -11456.56
$ws.Range("N86").Value = -14798

$ws.Range("H89").Value = 12568.404
$ws.Range("I89").Value = 12579.56
$ws.Range("J89").Value = 12552
$ws.Range("K89").Value = 62897.8
$ws.Range("L89").Value = 62760
$ws.Range("M89").Value = -57281.8
$ws.Range("N89").Value = -73992

$ws.Range("H105").Value = 2051.1428
$ws.Range("I105").Value = 2175.8
$ws.Range("J105").Value = 1739.5
$ws.Range("K105").Value = 2175.8
$ws.Range("L105").Value = 1739.5
$ws.Range("M105").Value = -428.8000000000002
$ws.Range("N105").Value = -5233.5

$ws.Range("H134").Value = 1099.25
$ws.Range("I134").Value = 1099.25
$ws.Range("K134").Value = 3297.75
$ws.Range("M134").Value = -762.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1341.5
$ws.Range("I31").Value = 1341.5
$ws.Range("K31").Value = 1341.5
$ws.Range("M31").Value = -1046.5

$ws.Range("H34").Value = 1341.5
$ws.Range("I34").Value = 1341.5
$ws.Range("K34").Value = 1341.5
$ws.Range("M34").Value = -1139.5

$ws.Range("H86").Value = 14833.167
$ws.Range("I86").Value = 13111.777
$ws.Range("K86").Value = 13111.777
$ws.Range("M86").Value = -11988.777

$ws.Range("H89").Value = 14833.167
$ws.Range("I89").Value = 13111.777
$ws.Range("K89").Value = 65558.88499999999
$ws.Range("M89").Value = -59942.88499999999

$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()

$ws.Range("H132").Value = 20012262
$ws.Range("I132").Value = 22235652
$ws.Range("J132").Value = 1750
$ws.Range("K132").Value = 66706956
$ws.Range("L132").Value = 5250
$ws.Range("M132").Value = -66704426
$ws.Range("N132").Value = -10310

$ws.Range("H134").Value = 1923.6364
$ws.Range("I134").Value = 1990.25
$ws.Range("J134").Value = 1746
$ws.Range("K134").Value = 5970.75
$ws.Range("L134").Value = 5238
$ws.Range("M134").Value = -3435.75
$ws.Range("N134").Value = -10308

$ws.Range("H141").Value = 347351.34
$ws.Range("J141").Value = 384520.25
$ws.Range("L141").Value = 384520.25
$ws.Range("N141").Value = -394880.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 37987.516
$ws.Range("J34").Value = 40762.555
$ws.Range("L34").Value = 122287.665
$ws.Range("N34").Value = -122455.665

$ws.Range("H106").Value = 17750

$ws.Range("H137").Value = 5472.8
$ws.Range("I137").Value = 3742.1428
$ws.Range("J137").Value = 6987.125
$ws.Range("K137").Value = 11226.4284
$ws.Range("L137").Value = 20961.375
$ws.Range("M137").Value = -6126.428400000001
$ws.Range("N137").Value = -31161.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 10000
$ws.Range("J49").Value = 10000
$ws.Range("L49").Value = 10000
$ws.Range("N49").Value = -10368

$ws.Range("H70").Value = 10000
$ws.Range("J70").Value = 10000
$ws.Range("L70").Value = 10000
$ws.Range("N70").Value = -10540

$ws.Range("H73").Value = 10000
$ws.Range("J73").Value = 10000
$ws.Range("L73").Value = 10000
$ws.Range("N73").Value = -11872

$ws.Range("H132").Value = 12824162
$ws.Range("I132").Value = 3649
$ws.Range("J132").Value = 41670316
$ws.Range("K132").Value = 10947
$ws.Range("L132").Value = 125010948
$ws.Range("M132").Value = -8417
$ws.Range("N132").Value = -125016008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 533.5
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H47").Value = 28500
$ws.Range("I47").Value = 7000
$ws.Range("J47").Value = 50000
$ws.Range("K47").Value = 7000
$ws.Range("L47").Value = 50000
$ws.Range("M47").Value = -6510
$ws.Range("N47").Value = -50980

$ws.Range("H52").Value = 28500
$ws.Range("I52").Value = 7000
$ws.Range("J52").Value = 50000
$ws.Range("K52").Value = 7000
$ws.Range("L52").Value = 50000
$ws.Range("M52").Value = -6767
$ws.Range("N52").Value = -50466

$ws.Range("H132").Value = 3186.0908
$ws.Range("I132").Value = 2966.1667
$ws.Range("J132").Value = 3268.5625
$ws.Range("K132").Value = 8898.500100000001
$ws.Range("L132").Value = 9805.6875
$ws.Range("M132").Value = -6368.500100000001
$ws.Range("N132").Value = -14865.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4338.2
$ws.Range("J107").Value = 3003
$ws.Range("L107").Value = 9009
$ws.Range("N107").Value = -12849

$ws.Range("H122").Value = 2712.4443
$ws.Range("I122").Value = 2716
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 8148
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -5698
$ws.Range("N122").Value = -13000

$ws.Range("H132").Value = 45465190
$ws.Range("I132").Value = 12032.842
$ws.Range("J132").Value = 333335200
$ws.Range("K132").Value = 36098.526
$ws.Range("L132").Value = 1000005600
$ws.Range("M132").Value = -33568.526
$ws.Range("N132").Value = -1000010660

$ws.Range("H136").Value = 6442.5713
$ws.Range("I136").Value = 7359.222
$ws.Range("K136").Value = 22077.666
$ws.Range("M136").Value = -19527.666
